$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 61, shifting all the
# subsequent rows (previously 61-143) down by one (now 62-144).
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row with its data.
$ws.Cells.Item(61, 1).Value = 8
$ws.Cells.Item(61, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(61, 3).Value = "Coquimbo"
$ws.Cells.Item(61, 4).Value = 44740
$ws.Cells.Item(61, 5).Value = 4
$ws.Cells.Item(61, 6).Value = 100112001
$ws.Cells.Item(61, 7).Value = "Berenjena"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 460
$ws.Cells.Item(61, 11).Value = 9000
$ws.Cells.Item(61, 12).Value = 10000
$ws.Cells.Item(61, 13).Value = 9500
$ws.Cells.Item(61, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(61, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(61, 16).Value = 190
$ws.Cells.Item(61, 17).Value = 50
$ws.Cells.Item(61, 18).Value = "Hortaliza"
